$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1017.6087
$ws.Range("I38").Value = 270.55
$ws.Range("K38").Value = 811.6500000000001
$ws.Range("M38").Value = -439.6500000000001
$ws.Range("H53").Value = 281.54544
$ws.Range("I53").Value = 350.875
$ws.Range("J53").Value = 96.666664
$ws.Range("K53").Value = 350.875
$ws.Range("L53").Value = 96.666664
$ws.Range("M53").Value = 286.125
$ws.Range("N53").Value = -1370.666664
$ws.Range("H64").Value = 4341.5625
$ws.Range("I64").Value = 3949.25
$ws.Range("K64").Value = 3949.25
$ws.Range("M64").Value = -3701.25
$ws.Range("H67").Value = 4341.5625
$ws.Range("I67").Value = 3949.25
$ws.Range("K67").Value = 3949.25
$ws.Range("M67").Value = -3091.25
$ws.Range("H108").Value = 59957.5
$ws.Range("J108").Value = 59957.5
$ws.Range("L108").Value = 59957.5
$ws.Range("N108").Value = -67637.5
$ws.Range("H111").Value = 1063.8572
$ws.Range("I111").Value = 1063.8572
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3191.5716
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -124.5715999999998
$ws.Range("N111").Value = ""
$ws.Range("H116").Value = 6625.4287
$ws.Range("I116").Value = 4641.1113
$ws.Range("K116").Value = 4641.1113
$ws.Range("M116").Value = -1199.1113
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H128").Value = 100000
$ws.Range("J128").Value = 100000
$ws.Range("L128").Value = 100000
$ws.Range("N128").Value = -109960
$ws.Range("H134").Value = 82666
$ws.Range("J134").Value = 82666
$ws.Range("L134").Value = 82666
$ws.Range("N134").Value = -92806
$ws.Range("H137").Value = 10198.167
$ws.Range("J137").Value = 14461.741
$ws.Range("L137").Value = 43385.223
$ws.Range("N137").Value = -48485.223
$ws.Range("H140").Value = 84442.664
$ws.Range("J140").Value = 84442.664
$ws.Range("L140").Value = 84442.664
$ws.Range("N140").Value = -94802.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 12463.846
$ws.Range("I45").Value = 13168
$ws.Range("K45").Value = 13168
$ws.Range("M45").Value = -12791
$ws.Range("H105").Value = 37863
$ws.Range("J105").Value = 37863
$ws.Range("L105").Value = 37863
$ws.Range("N105").Value = -44851
$ws.Range("H122").Value = 2058.65
$ws.Range("I122").Value = 1971.4117
$ws.Range("J122").Value = 2553
$ws.Range("K122").Value = 5914.2351
$ws.Range("L122").Value = 7659
$ws.Range("M122").Value = -3464.2351
$ws.Range("N122").Value = -12559

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 10041
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 10041
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 10041
$ws.Range("M36").Value = ""
$ws.Range("N36").Value = -11109
$ws.Range("H94").Value = 3279
$ws.Range("I94").Value = 2288.8333
$ws.Range("K94").Value = 2288.8333
$ws.Range("M94").Value = -1837.8333
$ws.Range("H103").Value = 34881
$ws.Range("J103").Value = 34881
$ws.Range("L103").Value = 34881
$ws.Range("N103").Value = -37225

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 54110.8
$ws.Range("J106").Value = 54110.8
$ws.Range("L106").Value = 54110.8
$ws.Range("N106").Value = -56634.8
$ws.Range("H107").Value = 1876.2727
$ws.Range("I107").Value = 1926.5714
$ws.Range("J107").Value = 1788.25
$ws.Range("K107").Value = 1926.5714
$ws.Range("L107").Value = 1788.25
$ws.Range("M107").Value = -6.57140000000004
$ws.Range("N107").Value = -5628.25
$ws.Range("H122").Value = 3498.182
$ws.Range("I122").Value = 2846.6
$ws.Range("K122").Value = 8539.8
$ws.Range("M122").Value = -6089.799999999999
$ws.Range("H141").Value = 372132.56
$ws.Range("J141").Value = 413322.34
$ws.Range("L141").Value = 413322.34
$ws.Range("N141").Value = -423682.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2764.7307
$ws.Range("J34").Value = 5770.0835
$ws.Range("L34").Value = 17310.2505
$ws.Range("N34").Value = -17478.2505
$ws.Range("H59").Value = 1466.6666
$ws.Range("J59").Value = 1500
$ws.Range("L59").Value = 4500
$ws.Range("N59").Value = -5580
$ws.Range("H63").Value = 2999
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 2998
$ws.Range("K63").Value = 9000
$ws.Range("L63").Value = 8994
$ws.Range("M63").Value = -8251
$ws.Range("N63").Value = -10492
$ws.Range("H66").Value = 2999
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 2998
$ws.Range("K66").Value = 27000
$ws.Range("L66").Value = 26982
$ws.Range("M66").Value = -23256
$ws.Range("N66").Value = -34470
$ws.Range("H107").Value = 1149.4375
$ws.Range("J107").Value = 496.14285
$ws.Range("L107").Value = 1488.42855
$ws.Range("N107").Value = -5328.428550000001
$ws.Range("H123").Value = 4000
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1060
$ws.Range("N126").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1078.625
$ws.Range("J97").Value = 995
$ws.Range("L97").Value = 995
$ws.Range("N97").Value = -1987
$ws.Range("H132").Value = 13195.407
$ws.Range("I132").Value = 9864.738
$ws.Range("K132").Value = 29594.214
$ws.Range("M132").Value = -27064.214
$ws.Range("H139").Value = 94278.86
$ws.Range("J139").Value = 94278.86
$ws.Range("L139").Value = 94278.86
$ws.Range("N139").Value = -104558.86
$ws.Range("H141").Value = 75844
$ws.Range("J141").Value = 75844
$ws.Range("L141").Value = 75844
$ws.Range("N141").Value = -86204

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2867.6667
$ws.Range("I7").Value = 2852
$ws.Range("J7").Value = 2899
$ws.Range("K7").Value = 2852
$ws.Range("L7").Value = 2899
$ws.Range("M7").Value = -2740
$ws.Range("N7").Value = -3123
$ws.Range("H22").Value = 1913.9546
$ws.Range("I22").Value = 1816
$ws.Range("J22").Value = 1942.7646
$ws.Range("K22").Value = 1816
$ws.Range("L22").Value = 1942.7646
$ws.Range("M22").Value = -1521
$ws.Range("N22").Value = -2532.7646
$ws.Range("H27").Value = 1913.9546
$ws.Range("I27").Value = 1816
$ws.Range("J27").Value = 1942.7646
$ws.Range("K27").Value = 1816
$ws.Range("L27").Value = 1942.7646
$ws.Range("M27").Value = -1709
$ws.Range("N27").Value = -2156.7646
$ws.Range("H74").Value = 27784.084
$ws.Range("I74").Value = 21774.125
$ws.Range("K74").Value = 21774.125
$ws.Range("M74").Value = -20776.125
$ws.Range("H77").Value = 27784.084
$ws.Range("I77").Value = 21774.125
$ws.Range("K77").Value = 65322.375
$ws.Range("M77").Value = -60330.375
$ws.Range("H122").Value = 4473.2
$ws.Range("I122").Value = 6241.6
$ws.Range("J122").Value = 3589
$ws.Range("K122").Value = 18724.8
$ws.Range("L122").Value = 10767
$ws.Range("M122").Value = -16274.8
$ws.Range("N122").Value = -15667
$ws.Range("H126").Value = 2867.6667
$ws.Range("I126").Value = 2852
$ws.Range("J126").Value = 2899
$ws.Range("K126").Value = 8556
$ws.Range("L126").Value = 8697
$ws.Range("M126").Value = -6086
$ws.Range("N126").Value = -13637
$ws.Range("H132").Value = 5076.421
$ws.Range("I132").Value = 5379.9614
$ws.Range("J132").Value = 4418.75
$ws.Range("K132").Value = 16139.8842
$ws.Range("L132").Value = 13256.25
$ws.Range("M132").Value = -13609.8842
$ws.Range("N132").Value = -18316.25
$ws.Range("H136").Value = 3941.7446
$ws.Range("I136").Value = 3539.0857
$ws.Range("K136").Value = 10617.2571
$ws.Range("M136").Value = -8067.257100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2514.7058
$ws.Range("I96").Value = 2559.6
$ws.Range("K96").Value = 2559.6
$ws.Range("M96").Value = -1186.6
$ws.Range("H107").Value = 5557350
$ws.Range("I107").Value = 1238.3572
$ws.Range("K107").Value = 3715.0716
$ws.Range("M107").Value = -1795.0716
$ws.Range("H124").Value = 44950
$ws.Range("J124").Value = 44950
$ws.Range("L124").Value = 44950
$ws.Range("N124").Value = -54770
$ws.Range("H132").Value = 175072.62
$ws.Range("I132").Value = 326023.3
$ws.Range("J132").Value = 24121.938
$ws.Range("K132").Value = 978069.8999999999
$ws.Range("L132").Value = 72365.814
$ws.Range("M132").Value = -975539.8999999999
$ws.Range("N132").Value = -77425.814
$ws.Range("H136").Value = 2747.1765
$ws.Range("I136").Value = 2436.5
$ws.Range("J136").Value = 2964.65
$ws.Range("K136").Value = 7309.5
$ws.Range("L136").Value = 8893.95
$ws.Range("M136").Value = -4759.5
$ws.Range("N136").Value = -13993.95
$ws.Range("H138").Value = 60000
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

Write-Output "Applied 237 cell updates across 8 sheets."